$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row for "하나스팩33호" above the current row 5 ("신한스팩13호"),
#    shifting all subsequent rows (previously 5-21) down to 6-22.
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "하나스팩33호"
$ws.Range("B5").Value = "2024.04.08~04.09"
$ws.Range("C5").Value = "2,000~2,000"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 7000
$ws.Range("F5").Value = "하나증권"

# 2) "하나스팩32호" (now row 9) had its 확정공모가 (column D) finalized from "-" to "2000".
#    Force the cell to stay text (matching the rest of column D), then restore the
#    default "Normal" style so no stray formatting is introduced.
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2000"
$ws.Range("D9").Style = "Normal"

# 3) The last data row ("코셈", now shifted to row 22) was removed entirely.
$ws.Rows.Item(22).Delete()
